# Scheduled runner update: refresh Universalis market-price snapshots
# (currentAveragePrice / LevePrice / LeveProfit columns) across all Leve sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 563.4375
$ws.Range("I19").Value = 591.4545000000001
$ws.Range("J19").Value = 501.8
$ws.Range("K19").Value = 591.4545000000001
$ws.Range("L19").Value = 501.8
$ws.Range("M19").Value = -416.4545000000001
$ws.Range("N19").Value = -851.8

$ws.Range("H88").Value = 806.875
$ws.Range("I88").Value = 437.6
$ws.Range("J88").Value = 974.7273
$ws.Range("K88").Value = 437.6
$ws.Range("L88").Value = 974.7273
$ws.Range("M88").Value = -31.60000000000002
$ws.Range("N88").Value = -1786.7273

$ws.Range("H91").Value = 806.875
$ws.Range("I91").Value = 437.6
$ws.Range("J91").Value = 974.7273
$ws.Range("K91").Value = 437.6
$ws.Range("L91").Value = 974.7273
$ws.Range("M91").Value = 966.4
$ws.Range("N91").Value = -3782.7273

$ws.Range("H92").Value = 3428.75
$ws.Range("I92").Value = 3643.111
$ws.Range("K92").Value = 3643.111
$ws.Range("M92").Value = -2395.111

$ws.Range("H96").Value = 994.1539
$ws.Range("J96").Value = 528.5
$ws.Range("L96").Value = 1585.5
$ws.Range("N96").Value = -4331.5

$ws.Range("H100").Value = 1614.75
$ws.Range("I100").Value = 1465
$ws.Range("J100").Value = 1839.375
$ws.Range("K100").Value = 1465
$ws.Range("L100").Value = 1839.375
$ws.Range("M100").Value = -924
$ws.Range("N100").Value = -2921.375

$ws.Range("H132").Value = 11616.4375
$ws.Range("I132").Value = 9289.704
$ws.Range("J132").Value = 24180.8
$ws.Range("K132").Value = 27869.112
$ws.Range("L132").Value = 72542.39999999999
$ws.Range("M132").Value = -25339.112
$ws.Range("N132").Value = -77602.39999999999

$ws.Range("H138").Value = 3912.4167
$ws.Range("I138").Value = 5193.8184
$ws.Range("K138").Value = 15581.4552
$ws.Range("M138").Value = -10441.4552

$ws.Range("H141").Value = 4665.4614
$ws.Range("I141").Value = 3969.652
$ws.Range("K141").Value = 11908.956
$ws.Range("M141").Value = -6728.956

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 53332.668
$ws.Range("I55").Value = 50000
$ws.Range("K55").Value = 50000
$ws.Range("M55").Value = -49685

$ws.Range("H61").Value = 11088.023
$ws.Range("I61").Value = 7266.1055
$ws.Range("J61").Value = 14245.261
$ws.Range("K61").Value = 7266.1055
$ws.Range("L61").Value = 14245.261
$ws.Range("M61").Value = -7054.1055
$ws.Range("N61").Value = -14669.261

$ws.Range("H97").Value = 2476.037
$ws.Range("I97").Value = 700.05554
$ws.Range("J97").Value = 6028
$ws.Range("K97").Value = 700.05554
$ws.Range("L97").Value = 6028
$ws.Range("M97").Value = -204.05554
$ws.Range("N97").Value = -7020

$ws.Range("H132").Value = 2048766.2
$ws.Range("I132").Value = 3106.4
$ws.Range("K132").Value = 9319.200000000001
$ws.Range("M132").Value = -6789.200000000001

$ws.Range("H136").Value = 11088.023
$ws.Range("I136").Value = 7266.1055
$ws.Range("J136").Value = 14245.261
$ws.Range("K136").Value = 21798.3165
$ws.Range("L136").Value = 42735.783
$ws.Range("M136").Value = -19248.3165
$ws.Range("N136").Value = -47835.783

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1307.7084
$ws.Range("J80").Value = 1498.65
$ws.Range("L80").Value = 1498.65
$ws.Range("N80").Value = -3494.65

$ws.Range("H83").Value = 1307.7084
$ws.Range("J83").Value = 1498.65
$ws.Range("L83").Value = 7493.25
$ws.Range("N83").Value = -17477.25

$ws.Range("H105").Value = 2464.8333
$ws.Range("I105").Value = 1957.8
$ws.Range("K105").Value = 1957.8
$ws.Range("M105").Value = -210.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19963
$ws.Range("J31").Value = 19979.4
$ws.Range("L31").Value = 19979.4
$ws.Range("N31").Value = -20569.4

$ws.Range("H34").Value = 19963
$ws.Range("J34").Value = 19979.4
$ws.Range("L34").Value = 19979.4
$ws.Range("N34").Value = -20383.4

$ws.Range("H104").Value = 20000
$ws.Range("J104").Value = 20000
$ws.Range("L104").Value = 20000
$ws.Range("N104").Value = -25242

$ws.Range("H105").Value = 16085.8
$ws.Range("I105").Value = 25664.75
$ws.Range("J105").Value = 9699.833000000001
$ws.Range("K105").Value = 25664.75
$ws.Range("L105").Value = 9699.833000000001
$ws.Range("M105").Value = -23917.75
$ws.Range("N105").Value = -13193.833

$ws.Range("H122").Value = 3984.8057
$ws.Range("I122").Value = 2258.7097
$ws.Range("K122").Value = 6776.1291
$ws.Range("M122").Value = -4326.1291

$ws.Range("H132").Value = 6836.5713
$ws.Range("I132").Value = 3078.8572
$ws.Range("K132").Value = 9236.571599999999
$ws.Range("M132").Value = -6706.571599999999

$ws.Range("H141").Value = 258196.75
$ws.Range("J141").Value = 258196.75
$ws.Range("L141").Value = 258196.75
$ws.Range("N141").Value = -268556.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 240.27777
$ws.Range("I14").Value = 240.27777
$ws.Range("K14").Value = 720.83331
$ws.Range("M14").Value = -547.83331

$ws.Range("H17").Value = 102.1
$ws.Range("I17").Value = 41
$ws.Range("J17").Value = 142.83333
$ws.Range("K17").Value = 123
$ws.Range("L17").Value = 428.49999
$ws.Range("M17").Value = 46
$ws.Range("N17").Value = -766.49999

$ws.Range("H75").Value = 850.11536
$ws.Range("I75").Value = 654
$ws.Range("K75").Value = 1962
$ws.Range("M75").Value = -964

$ws.Range("H78").Value = 850.11536
$ws.Range("I78").Value = 654
$ws.Range("K78").Value = 5886
$ws.Range("M78").Value = -894

$ws.Range("H131").Value = 1474.0935
$ws.Range("J131").Value = 1481.2476
$ws.Range("L131").Value = 4443.7428
$ws.Range("N131").Value = -14523.7428

$ws.Range("H140").Value = 2822
$ws.Range("I140").Value = 1453.3334
$ws.Range("K140").Value = 4360.0002
$ws.Range("M140").Value = 819.9997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 1002.2273
$ws.Range("I2").Value = 1241.4706
$ws.Range("J2").Value = 188.8
$ws.Range("K2").Value = 1241.4706
$ws.Range("L2").Value = 188.8
$ws.Range("M2").Value = -1128.4706
$ws.Range("N2").Value = -414.8

$ws.Range("H97").Value = 5588.467
$ws.Range("J97").Value = 8402.429
$ws.Range("L97").Value = 8402.429
$ws.Range("N97").Value = -9394.429

$ws.Range("H113").Value = 130242.43
$ws.Range("J113").Value = 2500
$ws.Range("L113").Value = 2500
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13663.6875
$ws.Range("I22").Value = 17432.857
$ws.Range("J22").Value = 10732.111
$ws.Range("K22").Value = 17432.857
$ws.Range("L22").Value = 10732.111
$ws.Range("M22").Value = -17137.857
$ws.Range("N22").Value = -11322.111

$ws.Range("H27").Value = 13663.6875
$ws.Range("I27").Value = 17432.857
$ws.Range("J27").Value = 10732.111
$ws.Range("K27").Value = 17432.857
$ws.Range("L27").Value = 10732.111
$ws.Range("M27").Value = -17325.857
$ws.Range("N27").Value = -10946.111

$ws.Range("H40").Value = 8328.15
$ws.Range("I40").Value = 5271
$ws.Range("J40").Value = 17499.6
$ws.Range("K40").Value = 5271
$ws.Range("L40").Value = 17499.6
$ws.Range("M40").Value = -5135
$ws.Range("N40").Value = -17771.6

$ws.Range("H46").Value = 1255249.4
$ws.Range("J46").Value = 8999
$ws.Range("L46").Value = 8999
$ws.Range("N46").Value = -9375

$ws.Range("H55").Value = 217.77777
$ws.Range("I55").Value = 162.5
$ws.Range("J55").Value = 262
$ws.Range("K55").Value = 162.5
$ws.Range("L55").Value = 262
$ws.Range("M55").Value = 10.5
$ws.Range("N55").Value = -608

$ws.Range("H68").Value = 5890.6665
$ws.Range("I68").Value = 2193
$ws.Range("J68").Value = 6630.2
$ws.Range("K68").Value = 2193
$ws.Range("L68").Value = 6630.2
$ws.Range("M68").Value = -1444
$ws.Range("N68").Value = -8128.2

$ws.Range("H71").Value = 5890.6665
$ws.Range("I71").Value = 2193
$ws.Range("J71").Value = 6630.2
$ws.Range("K71").Value = 10965
$ws.Range("L71").Value = 33151
$ws.Range("M71").Value = -7221
$ws.Range("N71").Value = -40639

$ws.Range("H122").Value = 7316.516
$ws.Range("I122").Value = 5445.579
$ws.Range("J122").Value = 10278.833
$ws.Range("K122").Value = 16336.737
$ws.Range("L122").Value = 30836.499
$ws.Range("M122").Value = -13886.737
$ws.Range("N122").Value = -35736.499

$ws.Range("H136").Value = 10002.164
$ws.Range("I136").Value = 9538
$ws.Range("J136").Value = 10698.409
$ws.Range("K136").Value = 28614
$ws.Range("L136").Value = 32095.227
$ws.Range("M136").Value = -26064
$ws.Range("N136").Value = -37195.227

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 2275000
$ws.Range("I75").Value = 30000
$ws.Range("J75").Value = 3023333.2
$ws.Range("K75").Value = 30000
$ws.Range("L75").Value = 3023333.2
$ws.Range("M75").Value = -29064
$ws.Range("N75").Value = -3025205.2

$ws.Range("H78").Value = 2275000
$ws.Range("I78").Value = 30000
$ws.Range("J78").Value = 3023333.2
$ws.Range("K78").Value = 90000
$ws.Range("L78").Value = 9069999.600000001
$ws.Range("M78").Value = -85320
$ws.Range("N78").Value = -9079359.600000001

$ws.Range("H122").Value = 23259798
$ws.Range("I122").Value = 34484630
$ws.Range("J122").Value = 8355.929
$ws.Range("K122").Value = 103453890
$ws.Range("L122").Value = 25067.787
$ws.Range("M122").Value = -103451440
$ws.Range("N122").Value = -29967.787

$ws.Range("H132").Value = 11136.704
$ws.Range("I132").Value = 4967.1665
$ws.Range("K132").Value = 14901.4995
$ws.Range("M132").Value = -12371.4995
